# Update the "Sprint Number" values (column F) on the Product Backlog sheet
# and give the whole Sprint Number column an integer ("0") number format,
# matching the author's edit to docs/Product Backlog-1.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Sprint Number (column F) value actually changed.
$updates = @{
    7  = 1
    10 = 2
    14 = 2
    19 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 1
    26 = 1
    29 = 1
    30 = 1
    31 = 1
    33 = 1
    34 = 1
    35 = 1
    45 = 2
    48 = 2
    49 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}

# Apply an integer number format to the entire Sprint Number data range,
# which causes Excel to create/assign a dedicated style (numFmtId 1) for
# these cells, same as in the source edit.
$ws.Range("F4:F49").NumberFormat = "0"
